$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the "Arquivo" (file) column
$ws.Cells.Item(1, 4).Value = "Arquivo"

# Update row 2 data
$ws.Cells.Item(2, 1).Value = "Isa"
$ws.Cells.Item(2, 2).Value = "5511943808142"
$ws.Cells.Item(2, 3).Value = "S"

# Update row 3 data
$ws.Cells.Item(3, 1).Value = "Pietra"
$ws.Cells.Item(3, 2).Value = "5511941900392"
$ws.Cells.Item(3, 3).Value = "S"

# Row 4 stays the same except new "Arquivo" value
$ws.Cells.Item(4, 4).Value = "teste.jpg"

$ws.Range("D5").Select() | Out-Null
